$d = $word.ActiveDocument

# Locate the paragraph that currently ends the "Dia 17/09..." entry so we can
# append the new "Dia 18/09..." metrics line right after it, before the
# trailing blank paragraph / sectPr.
$targetText = "Dia 17/09: 2hr (1 dia)"
$count = $d.Paragraphs.Count
$targetIdx = -1

for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text.TrimEnd([char]13, [char]7)
    if ($t -eq $targetText) {
        $targetIdx = $i
    }
}

if ($targetIdx -eq -1) {
    throw "Could not find paragraph with text '$targetText'"
}

$target = $d.Paragraphs.Item($targetIdx)

# Insert a new paragraph right after the target; it inherits the target's
# paragraph formatting (spacing, justification, Arial 12pt run formatting).
$target.Range.InsertParagraphAfter()

$newPara = $d.Paragraphs.Item($targetIdx + 1)
$newPara.Range.Text = "Dia 18/09: 1hr (1 dia)"
